# Append four new daily-report rows (180-183) to the Arequipa COVID data
# sheet, continuing straight on from the existing last row (179).
#
# Column layout (A:Q):
#   date, total_muestras, casos_positivos, casos_negativos, espera_resultado,
#   defunciones, recuperados, hospitalizados_positivos,
#   hospitalizados_sospechosos, uci_positivos, uci_sospechosos, trauma_shock,
#   defunciones_minsa, defunciones_essalud, defunciones_clinicas,
#   defunciones_domicilio, defunciones_sanidades

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 180 and 182/183 share the same look as row 179 (date col right-aligned
# date format, data cols in the "plain" style). Row 181 shares the look of
# row 178 (data cols in the "right aligned" style). Clone formatting first so
# new cells pick up the right style indices, then fill in the values.
$ws.Range("A179:Q179").Copy()
$ws.Range("A180:Q180").PasteSpecial($xlPasteFormats)
$ws.Range("A179:Q179").Copy()
$ws.Range("A182:Q183").PasteSpecial($xlPasteFormats)

$ws.Range("A178:Q178").Copy()
$ws.Range("A181:Q181").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

$data = @{
    180 = @(44130, 875291, 136263, 738764, 264, 2225, 128680, 227, 67, 47, 0, 3, 839, 1182, 39, 123, 42)
    181 = @(44131, 878156, 136488, 741370, 298, 2225, 129135, 223, 63, 49, 0, 3, 839, 1182, 39, 123, 42)
    182 = @(44132, 882024, 136776, 744993, 255, 2229, 129896, 228, 57, 48, 0, 3, 841, 1184, 39, 123, 42)
    183 = @(44133, 885094, 136965, 747906, 223, 2231, 130198, 231, 28, 51, 0, 3, 842, 1185, 39, 123, 42)
}

foreach ($r in 180, 181, 182, 183) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

Write-Output "Added rows 180-183"
